$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 940.5
$ws.Range("I6").Value = 464
$ws.Range("J6").Value = 2052.3333
$ws.Range("K6").Value = 1392
$ws.Range("L6").Value = 6156.999899999999
$ws.Range("M6").Value = -1280
$ws.Range("N6").Value = -6380.999899999999

# Row 86
$ws.Range("H86").Value = 29750
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 29750
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 29750
$ws.Range("N86").Value = -31996

# Row 89
$ws.Range("H89").Value = 29750
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 29750
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 148750
$ws.Range("N89").Value = -159982

# Row 92
$ws.Range("H92").Value = 1343.5454
$ws.Range("I92").Value = 1149.8334
$ws.Range("J92").Value = 1576
$ws.Range("K92").Value = 1149.8334
$ws.Range("L92").Value = 1576
$ws.Range("M92").Value = 98.16660000000002
$ws.Range("N92").Value = -4072

# Row 100
$ws.Range("H100").Value = 681.125
$ws.Range("I100").Value = 681.125
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 681.125
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -140.125
$ws.Range("N100").ClearContents()

# Row 135
$ws.Range("H135").Value = 1036
$ws.Range("I135").Value = 795
$ws.Range("J135").Value = 2000
$ws.Range("K135").Value = 7155
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = -4620

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 978.6667
$ws.Range("I32").Value = 917.4
$ws.Range("J32").Value = 1285
$ws.Range("K32").Value = 917.4
$ws.Range("L32").Value = 1285
$ws.Range("M32").Value = -630.4
$ws.Range("N32").Value = -1859

# Row 61
$ws.Range("H61").Value = 418.3
$ws.Range("I61").Value = 418.3
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 418.3
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -206.3

# Row 136
$ws.Range("H136").Value = 418.3
$ws.Range("I136").Value = 418.3
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 1254.9
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 1295.1

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 1585.8572
$ws.Range("I20").Value = 1585.8572
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 1585.8572
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -1338.8572
$ws.Range("N20").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1389.6923
$ws.Range("I31").Value = 1119.5555
$ws.Range("J31").Value = 1997.5
$ws.Range("K31").Value = 1119.5555
$ws.Range("L31").Value = 1997.5
$ws.Range("M31").Value = -824.5554999999999
$ws.Range("N31").Value = -2587.5

# Row 34
$ws.Range("H34").Value = 1389.6923
$ws.Range("I34").Value = 1119.5555
$ws.Range("J34").Value = 1997.5
$ws.Range("K34").Value = 1119.5555
$ws.Range("L34").Value = 1997.5
$ws.Range("M34").Value = -917.5554999999999
$ws.Range("N34").Value = -2401.5

# Row 58
$ws.Range("H58").Value = 2272.25
$ws.Range("I58").Value = 2168.2856
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 2168.2856
$ws.Range("L58").Value = 3000
$ws.Range("M58").Value = -1965.2856

# Row 132
$ws.Range("H132").Value = 1410.6
$ws.Range("I132").Value = 1300.75
$ws.Range("J132").Value = 1850
$ws.Range("K132").Value = 3902.25
$ws.Range("L132").Value = 5550
$ws.Range("M132").Value = -1372.25
$ws.Range("N132").Value = -10610

# Row 134
$ws.Range("H134").Value = 1753.2
$ws.Range("I134").Value = 1753.2
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 5259.6
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -2724.6

# Row 136
$ws.Range("H136").Value = 2272.25
$ws.Range("I136").Value = 2168.2856
$ws.Range("J136").Value = 3000
$ws.Range("K136").Value = 6504.8568
$ws.Range("L136").Value = 9000
$ws.Range("M136").Value = -3954.8568

$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 700
$ws.Range("I132").Value = 700
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2100
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 430

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 3248.8333
$ws.Range("I61").Value = 3248.8333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3248.8333
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3046.8333

# Row 64
$ws.Range("H64").Value = 45000
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 45000
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 45000
$ws.Range("N64").Value = -45450

# Row 67
$ws.Range("H67").Value = 45000
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 45000
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 45000
$ws.Range("N67").Value = -46560

# Row 68
$ws.Range("H68").Value = 1640.4
$ws.Range("I68").Value = 1800.5
$ws.Range("J68").Value = 1000
$ws.Range("K68").Value = 1800.5
$ws.Range("L68").Value = 1000
$ws.Range("M68").Value = -1051.5
$ws.Range("N68").Value = -2498

# Row 71
$ws.Range("H71").Value = 1640.4
$ws.Range("I71").Value = 1800.5
$ws.Range("J71").Value = 1000
$ws.Range("K71").Value = 9002.5
$ws.Range("L71").Value = 5000
$ws.Range("M71").Value = -5258.5
$ws.Range("N71").Value = -12488

# Row 100
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

# Row 113
$ws.Range("H113").Value = 3248.8333
$ws.Range("I113").Value = 3248.8333
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3248.8333
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1078.8333

# Row 132
$ws.Range("H132").Value = 3271.3333
$ws.Range("I132").Value = 3004
$ws.Range("J132").Value = 3405
$ws.Range("K132").Value = 9012
$ws.Range("L132").Value = 10215
$ws.Range("M132").Value = -6482
$ws.Range("N132").Value = -15275

# Row 136
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -12450

$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 574.75
$ws.Range("I96").Value = 224.5
$ws.Range("J96").Value = 925
$ws.Range("K96").Value = 224.5
$ws.Range("L96").Value = 925
$ws.Range("M96").Value = 1148.5
$ws.Range("N96").Value = -3671

# Row 100
$ws.Range("H100").Value = 215.66667
$ws.Range("I100").Value = 223.5
$ws.Range("J100").Value = 200
$ws.Range("K100").Value = 447
$ws.Range("L100").Value = 400
$ws.Range("M100").Value = 94

# Row 136
$ws.Range("H136").Value = 2800
$ws.Range("I136").Value = 2800
$ws.Range("J136").Value = 2800
$ws.Range("K136").Value = 8400
$ws.Range("L136").Value = 8400
$ws.Range("M136").Value = -5850
$ws.Range("N136").Value = -13500
